# Fill in the "selected IO / configured function / label" columns
# for the rows that correspond to buzzer, LEDs, buttons and other
# peripherals added as part of the OLED showcase wiring.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("IO分配情况")

# Row 4: PA13 - soft I2C reset line
$ws.Range("D4").Value = "PA13"
$ws.Range("E4").Value = "OUTPUT"
$ws.Range("F4").Value = "RST_soft_I2C"

# Row 5: PA14 - high voltage power supply control
$ws.Range("D5").Value = "PA14"
$ws.Range("E5").Value = "OUTPUT"
$ws.Range("F5").Value = "HV_PSU"

# Row 7: PA8 - Geiger-Muller tube pulse input
$ws.Range("D7").Value = "PA8"
$ws.Range("E7").Value = "INPUT"
$ws.Range("F7").Value = "GM_pulse"

# Row 10: PB2 - power button (wakeup)
$ws.Range("D10").Value = "PB2"
$ws.Range("E10").Value = "WKUP_WKUP2"
$ws.Range("F10").Value = "POWER_button"

# Row 11: PB3 - power enable
$ws.Range("D11").Value = "PB3"
$ws.Range("E11").Value = "OUTPUT"
$ws.Range("F11").Value = "POWER_EN"

# Row 12: PB8 - USB sense
$ws.Range("D12").Value = "PB8"
$ws.Range("E12").Value = "INPUT"
$ws.Range("F12").Value = "USB_sense"

# Row 15: PB11 - buzzer
$ws.Range("D15").Value = "PB11"
$ws.Range("E15").Value = "OUTPUT"
$ws.Range("F15").Value = "BUZZER"

# Row 22: PC8 - battery ADC input
$ws.Range("D22").Value = "PC8"
$ws.Range("E22").Value = "ADC_IN9"
$ws.Range("F22").Value = "BATT_ADC_IN9"

# Row 24: PC10 - menu button
$ws.Range("D24").Value = "PC10"
$ws.Range("E24").Value = "INPUT"
$ws.Range("F24").Value = "MENU_button"

# Row 25: label for PD9 (red LED)
$ws.Range("F25").Value = "LED1_red"

# Row 26: label for PD10 (green LED)
$ws.Range("F26").Value = "LED2_green"

# Row 31: label for PD0 (soft I2C SDA)
$ws.Range("F31").Value = "SDA_soft_I2C"

# Row 32: label for PD1 (soft I2C SCL)
$ws.Range("F32").Value = "SCL_soft_I2C"
